$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 526, shifting existing rows 526:601 down to 527:602.
$ws.Rows.Item(526).Insert()

# Populate the newly inserted row 526 with the new data record.
$ws.Range("A526").Value = 3
$ws.Range("B526").Value = "Femacal de La Calera"
$ws.Range("C526").Value = "Coquimbo"
$ws.Range("D526").Value = 45077
$ws.Range("E526").Value = 5
$ws.Range("F526").Value = 100112031
$ws.Range("G526").Value = "Poroto verde"
$ws.Range("H526").Value = "Magnum"
$ws.Range("I526").Value = "Primera"
$ws.Range("J526").Value = 78
$ws.Range("K526").Value = 34000
$ws.Range("L526").Value = 35000
$ws.Range("M526").Value = 34513
$ws.Range("N526").Value = "`$/saco 25 kilos"
$ws.Range("O526").Value = "Provincia de Limarí"
$ws.Range("P526").Value = 1381
$ws.Range("Q526").Value = 25
$ws.Range("R526").Value = "Hortaliza"
